# Completed SO fields additions
# Adds five new header columns (Deliver To Name, Product Code, Product Name,
# Order Load, Order Quantity) before the trailing "Remarks" column on the
# "customermaster" sheet, pushing Remarks from O1 to T1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "Remarks" out to the new last column, and fill the freed-up /
# newly-added header cells with the new field names (in order, so the
# shared-string table gets them appended in the same sequence as the
# original authoring tool produced).
$ws.Range("T1").Value = "Remarks"
$ws.Range("O1").Value = "Deliver To Name"
$ws.Range("P1").Value = "Product Code"
$ws.Range("Q1").Value = "Product Name"
$ws.Range("R1").Value = "Order Load"
$ws.Range("S1").Value = "Order Quantity"

# Match the bold header formatting used by every other column on row 1.
$ws.Range("O1:T1").Font.Bold = $true

# Best-fit-ish column widths for the new columns (character-width units;
# the stored OOXML width includes the standard padding offset).
$ws.Columns.Item(15).ColumnWidth = 15.0221354166667   # O - Deliver To Name
$ws.Columns.Item(16).ColumnWidth = 12.0221354166667   # P - Product Code
$ws.Columns.Item(17).ColumnWidth = 12.8776041666667   # Q - Product Name
$ws.Columns.Item(18).ColumnWidth = 9.87760416666667   # R - Order Load
$ws.Columns.Item(19).ColumnWidth = 13.5924479166667   # S - Order Quantity
$ws.Columns.Item(20).ColumnWidth = 7.73697916666667   # T - Remarks

# Restore the active selection to match the post-edit workbook state.
$ws.Range("K9").Select() | Out-Null
